# Apply cryptos.xlsx data refresh (Fri Jan 12 11:54:43 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.092.96"
$ws.Range("E2").Value = "  -1.96%  "
$ws.Range("D3").Value = "2.649.32"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  +0.07%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "310.72"
$c.ClearFormats()
$ws.Range("E5").Value = "  -0.99%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "98.73"
$c.ClearFormats()
$ws.Range("E6").Value = "  -5.36%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.599"
$c.ClearFormats()
$ws.Range("E7").Value = "  -2.00%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("E9").Value = "  -2.75%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "38.75"
$c.ClearFormats()
$ws.Range("E10").Value = "  -2.12%  "
$ws.Range("E11").Value = "  -1.15%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "54.29"
$c.ClearFormats()
$ws.Range("E12").Value = "  -1.76%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "8.10"
$c.ClearFormats()
$ws.Range("E13").Value = "  -3.28%  "
$ws.Range("D14").Value = "3.055.77"
$ws.Range("E14").Value = "  -0.07%  "
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("D16").Value = "2.646.65"
$ws.Range("E16").Value = "  -0.66%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.926"
$c.ClearFormats()
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "14.97"
$c.ClearFormats()
$ws.Range("E18").Value = "  -2.15%  "
$ws.Range("D19").Value = "46.120.51"
$ws.Range("E19").Value = "  -2.99%  "
$ws.Range("E20").Value = "  -1.48%  "
$ws.Range("E21").Value = "  -0.78%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "12.90"
$c.ClearFormats()
$ws.Range("E22").Value = "  -3.03%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "74.67"
$c.ClearFormats()
$ws.Range("E23").Value = "  +2.46%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "282.70"
$c.ClearFormats()
$ws.Range("E24").Value = "  +6.10%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "3.06"
$c.ClearFormats()
$ws.Range("E25").Value = "  -1.80%  "
$ws.Range("E26").Value = "  +0.36%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "30.42"
$c.ClearFormats()
$ws.Range("E27").Value = "  -4.41%  "
$ws.Range("E28").Value = "  -0.24%  "
$ws.Range("E29").Value = "  -1.71%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "38.71"
$c.ClearFormats()
$ws.Range("E30").Value = "  -6.64%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "2.25"
$c.ClearFormats()
$ws.Range("E31").Value = "  -2.99%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "6.28"
$c.ClearFormats()
$ws.Range("E32").Value = "  +0.18%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.74"
$c.ClearFormats()
$ws.Range("E33").Value = "  -1.80%  "
$ws.Range("E34").Value = "  +0.32%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "156.53"
$c.ClearFormats()
$ws.Range("E35").Value = "  +2.19%  "
$ws.Range("E37").Value = "  -1.51%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.123"
$c.ClearFormats()
$ws.Range("E38").Value = "  +2.71%  "
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "25.87"
$c.ClearFormats()
$ws.Range("E39").Value = "  +13.85%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.124"
$c.ClearFormats()
$ws.Range("E40").Value = "  -0.38%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "15.85"
$c.ClearFormats()
$ws.Range("E41").Value = "  -6.90%  "
$ws.Range("E42").Value = "  -1.06%  "
$ws.Range("E43").Value = "  -4.74%  "
$ws.Range("E44").Value = "  -6.75%  "
$ws.Range("D45").Value = "2.155.68"
$ws.Range("E45").Value = "  +2.87%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.ClearFormats()
$ws.Range("E46").Value = "  -0.08%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "94.05"
$c.ClearFormats()
$ws.Range("E47").Value = "  -0.77%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "9.30"
$c.ClearFormats()
$ws.Range("E48").Value = "  -0.05%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "111.29"
$c.ClearFormats()
$ws.Range("E49").Value = "  -4.03%  "
$ws.Range("D50").Value = "2.906.39"
$ws.Range("E50").Value = "  -0.28%  "
$ws.Range("E51").Value = "  -1.31%  "
